$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> EmployeeLoginCredentials ---
$wsLogin = $wb.Worksheets.Item(1)
$wsLogin.Name = "EmployeeLoginCredentials"

# --- Insert two new sheets after it: Login, AddReport ---
$wsLoginSheet = $wb.Worksheets.Add($null, $wsLogin)
$wsLoginSheet.Name = "Login"

$wsAddReport = $wb.Worksheets.Add($null, $wsLoginSheet)
$wsAddReport.Name = "AddReport"

# --- AddEmployee sheet is now the 4th tab ---
$wsAddEmployee = $wb.Worksheets.Item("AddEmployee")

# ============ EmployeeLoginCredentials (was Sheet1) ============
$wsLogin.Range("A1").Value = "FirstName"
$wsLogin.Range("B1").Value = "LastName"
$wsLogin.Range("C1").Value = "Username"
$wsLogin.Range("D1").Value = "Password"

$wsLogin.Range("A2").Value = "Noraxm"
$wsLogin.Range("B2").Value = "Felixxm"
$wsLogin.Range("C2").Value = "Adminxm"
$wsLogin.Range("D2").Value = "admin123_@H"

$wsLogin.Range("A3").Value = "Naomixm"
$wsLogin.Range("B3").Value = "Heightxm"
$wsLogin.Range("C3").Value = "Adminxm"
$wsLogin.Range("D3").Value = "admin123_@H"

$wsLogin.Range("A4").Value = "Mayaxm"
$wsLogin.Range("B4").Value = "Faradayxm"
$wsLogin.Range("C4").Value = "Royal1235xm"
$wsLogin.Range("D4").Value = "admin123_@H"

# Column widths (approx, matches autofit look)
$wsLogin.Columns.Item(1).ColumnWidth = 23.6
$wsLogin.Columns.Item(2).ColumnWidth = 26.5
$wsLogin.Columns.Item(3).ColumnWidth = 24.4
$wsLogin.Columns.Item(4).ColumnWidth = 26.7
$wsLogin.Columns.Item(5).ColumnWidth = 22.6

$wsLogin.Range("A1:D1").Font.Bold = $false
$wsLogin.Range("A1:D1").Font.Name = "Tahoma"
$wsLogin.Range("A1:D1").Font.Size = 16

$wsLogin.Range("E4").Select()

# ============ Login sheet (new) ============
$wsLoginSheet.Range("A1").Value = "Report Name"
$wsLoginSheet.Range("B1").Value = "Selection Criteria"
$wsLoginSheet.Range("C1").Value = "Selected Criteria Include"
$wsLoginSheet.Range("D1").Value = "Fields"

$wsLoginSheet.Columns.Item(1).ColumnWidth = 14.11
$wsLoginSheet.Columns.Item(2).ColumnWidth = 18.22
$wsLoginSheet.Columns.Item(3).ColumnWidth = 21.89
$wsLoginSheet.Columns.Item(4).ColumnWidth = 18.89
$wsLoginSheet.Columns.Item(5).ColumnWidth = 17.78

$wsLoginSheet.Range("D1").Select()

# ============ AddReport sheet (new, empty) ============
# left empty intentionally

# ============ AddEmployee (was sheet2) ============
$wsAddEmployee.Range("A1").Value = "FirstName"
$wsAddEmployee.Range("B1").Value = "MiddleName"
$wsAddEmployee.Range("C1").Value = "LastName"

$wsAddEmployee.Range("A2").Value = "Johnnn"
$wsAddEmployee.Range("B2").Value = "Mxx"
$wsAddEmployee.Range("C2").Value = "Smithxx"

$wsAddEmployee.Range("A3").Value = "Janeee"
$wsAddEmployee.Range("B3").Value = "Hxx"
$wsAddEmployee.Range("C3").Value = "Smithxx"

$wsAddEmployee.Columns.Item(2).ColumnWidth = 21.11

# Activate the Login tab as the selected sheet
$wsLoginSheet.Activate()
